# Apply "Updated symbol list" edits to the cryptos worksheet.
#
# All Price (column D) values are stored as text in the workbook (e.g. "235.70"),
# so before writing numeric-looking strings we force the cell's number format to
# Text ("@") to stop Excel from re-interpreting them as numbers and mangling the
# formatting (trailing zeros, rounding, scientific notation, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
}

# --- Price (column D) updates -------------------------------------------------
Set-TextValue "D2"  "235.70"
Set-TextValue "D3"  "21.72"
Set-TextValue "D7"  "6.459"
Set-TextValue "D8"  "0.8024"
Set-TextValue "D9"  "1.041"
Set-TextValue "D11" "0.07209"
Set-TextValue "D12" "0.03196"
Set-TextValue "D13" "0.02937"
Set-TextValue "D14" "0.09242"
Set-TextValue "D15" "0.001663"
Set-TextValue "D16" "3.256"
Set-TextValue "D17" "0.04781"
Set-TextValue "D20" "0.005075"
Set-TextValue "D21" "0.001049"

# --- Row 24 (LEO) ---------------------------------------------------------------
Set-TextValue "D24" "3.918"
$ws.Range("E24").Value = "23LEOLEOBestin24h"

# --- Rows 40-41 -------------------------------------------------------------
Set-TextValue "D40" "0.04114"
Set-TextValue "D41" "0.006997"

# --- Rows 42-43 swap (CEJI <-> BKEXToken) ------------------------------------
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1038"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002902"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining price updates -------------------------------------------------
Set-TextValue "D44" "0.008930"
Set-TextValue "D45" "0.00005439"
Set-TextValue "D48" "0.03324"
Set-TextValue "D50" "0.01011"
